$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Component header
$ws.Range("B2").Value = "Component: SCC (e-Bridge Cloud Client)"

# Row 6
$ws.Range("C6").Value = "MFP is powered on and in normal boot-up mode"
$ws.Range("D6").Value = "Enabling SCC function"
$ws.Range("E6").Value = "1. Access Self-diagnostic boot-up mode or service UI on panel`n2. Set 08-3820 from 0 (Disable) to 1 (Enable)`n3. Reboot the MFP"
$ws.Range("F6").Value = "SCC function should be enabled and MFP should start SCC process after reboot"

# Row 7
$ws.Range("C7").Value = "SCC function is enabled (08-3820 = 1)"
$ws.Range("D7").Value = "SCC Installation Report printing"
$ws.Range("E7").Value = "1. Change 08-3820 from 0 (Disable) to 1 (Enable)`n2. Reboot the MFP`n3. Wait for MFP registration processing to complete"
$ws.Range("F7").Value = "Installation report should be automatically printed with correct MFP information including Serial Number, MAC Address, and Registration Status"

# Row 8
$ws.Range("C8").Value = "SCC function is enabled"
$ws.Range("D8").Value = "SCC report printing language"
$ws.Range("E8").Value = "1. Change device language setting`n2. Enable SCC function and trigger installation report"
$ws.Range("F8").Value = "Installation report should be printed in English regardless of device language setting"

# Row 9
$ws.Range("C9").Value = "SCC function is enabled"
$ws.Range("D9").Value = "SCC report in job log"
$ws.Range("E9").Value = "1. Enable SCC function and trigger installation report`n2. Check job log"
$ws.Range("F9").Value = "SCC report printing should be recorded in Job Log as network print and user counter for built-in Admin should be incremented"

# Row 10
$ws.Range("C10").Value = "SCC function is enabled, network print restriction mode (08-9344) set to `"only private`""
$ws.Range("D10").Value = "SCC report with print restrictions"
$ws.Range("E10").Value = "1. Set 08-9344 to 1 (only private)`n2. Enable SCC function and trigger installation report"
$ws.Range("F10").Value = "SCC report should not be printed and job log should show error code 4221"

# Row 11
$ws.Range("C11").Value = "SCC function is enabled"
$ws.Range("D11").Value = "HTTPS communication protocol"
$ws.Range("E11").Value = "1. Enable SCC function`n2. Monitor network traffic during SCC communication`n3. Verify SSL/TLS connection to server"
$ws.Range("F11").Value = "MFP should use HTTPS for all communications with the SCC server"

# Row 12
$ws.Range("C12").Value = "SCC function is enabled"
$ws.Range("D12").Value = "Communication initiation"
$ws.Range("E12").Value = "1. Enable SCC function`n2. Monitor network traffic`n3. Observe communication patterns"
$ws.Range("F12").Value = "All communications should be initiated from the MFP side, never from the server side"

# Row 13
$ws.Range("C13").Value = "SCC function is enabled, proxy server required"
$ws.Range("D13").Value = "Proxy server authentication"
$ws.Range("E13").Value = "1. Configure proxy settings (08-3822 to 08-3826)`n2. Set proxy authentication credentials`n3. Trigger SCC communication"
$ws.Range("F13").Value = "MFP should successfully communicate with SCC server through the proxy with authentication"

# Row 14
$ws.Range("C14").Value = "SCC function is enabled, proxy server required without authentication"
$ws.Range("D14").Value = "Proxy server without authentication"
$ws.Range("E14").Value = "1. Configure proxy settings (08-3822 to 08-3826)`n2. Set proxy account ID to NULL`n3. Trigger SCC communication"
$ws.Range("F14").Value = "MFP should communicate with SCC server through proxy without authentication"

# Row 15
$ws.Range("C15").Value = "SCC function is enabled"
$ws.Range("D15").Value = "Port number fallback"
$ws.Range("E15").Value = "1. Block port 443 at firewall`n2. Allow port 8443`n3. Trigger SCC communication"
$ws.Range("F15").Value = "MFP should automatically use port 8443 when port 443 is unavailable"

# Row 16
$ws.Range("C16").Value = "MFP is in special startup mode"
$ws.Range("D16").Value = "SCC process startup"
$ws.Range("E16").Value = "1. Boot MFP in special startup mode`n2. Check if SCC process starts"
$ws.Range("F16").Value = "SCC process should not start in special startup mode"

# Row 17
$ws.Range("C17").Value = "SCC function is enabled"
$ws.Range("D17").Value = "Regular communication schedule - daily"
$ws.Range("E17").Value = "1. Set regular communication schedule to daily at specific time`n2. Wait for scheduled time`n3. Monitor network traffic"
$ws.Range("F17").Value = "MFP should initiate communication with SCC server at the specified time"

# Row 18
$ws.Range("C18").Value = "SCC function is enabled"
$ws.Range("D18").Value = "Regular communication schedule - interval"
$ws.Range("E18").Value = "1. Set regular communication schedule to interval (e.g., 15 min)`n2. Monitor network traffic"
$ws.Range("F18").Value = "MFP should initiate communication with SCC server at the specified interval"

# Row 19
$ws.Range("C19").Value = "SCC function is enabled"
$ws.Range("D19").Value = "Register Device - first time"
$ws.Range("E19").Value = "1. Enable SCC function on unregistered MFP`n2. Monitor communication with server"
$ws.Range("F19").Value = "MFP should send registration request with serial number and other required parameters, and receive token from server"

# Row 20
$ws.Range("C20").Value = "SCC function is enabled, MFP already registered"
$ws.Range("D20").Value = "Register Device - authentication"
$ws.Range("E20").Value = "1. Trigger communication on previously registered MFP`n2. Monitor communication with server"
$ws.Range("F20").Value = "MFP should authenticate using previously received token"
$ws.Range("H20").Value = ""

# Row 21
$ws.Range("C21").Value = "SCC function is enabled"
$ws.Range("D21").Value = "Server BUSY handling"
$ws.Range("E21").Value = "1. Simulate server BUSY response`n2. Monitor MFP behavior"
$ws.Range("F21").Value = "MFP should retry connection after waiting for the period specified by server, up to 3 retries"

# Row 22
$ws.Range("C22").Value = "SCC function is enabled"
$ws.Range("D22").Value = "Check for Updates"
$ws.Range("E22").Value = "1. Trigger regular communication`n2. Monitor Check for Updates request"
$ws.Range("F22").Value = "MFP should send firmware version and other parameters to check for available updates"

# Row 23
$ws.Range("C23").Value = "SCC function is enabled, update available"
$ws.Range("D23").Value = "Download Package - firmware"
$ws.Range("E23").Value = "1. Configure server to provide firmware update`n2. Trigger regular communication`n3. Monitor download process"
$ws.Range("F23").Value = "MFP should download update package and verify hash value"

# Row 24
$ws.Range("C24").Value = "SCC function is enabled, update available"
$ws.Range("D24").Value = "Install Package - firmware"
$ws.Range("E24").Value = "1. Download firmware update package`n2. Monitor installation process"
$ws.Range("F24").Value = "MFP should install firmware update at scheduled time and reboot"

# Row 25
$ws.Range("C25").Value = "SCC function is enabled, update available"
$ws.Range("D25").Value = "Update Status notification"
$ws.Range("E25").Value = "1. Install update package`n2. Monitor Update Status communication"
$ws.Range("F25").Value = "MFP should notify server of installation status with correct state transitions"

# Row 26
$ws.Range("C26").Value = "SCC function is enabled"
$ws.Range("D26").Value = "Send Baseline Data"
$ws.Range("E26").Value = "1. Trigger regular communication`n2. Monitor Send Baseline Data process"
$ws.Range("F26").Value = "MFP should collect and send baseline data according to SendDataConfig settings"

# Row 27
$ws.Range("C27").Value = "SCC function is enabled"
$ws.Range("D27").Value = "Send Regular Data"
$ws.Range("E27").Value = "1. Complete update installation`n2. Monitor Send Regular Data process"
$ws.Range("F27").Value = "MFP should collect and send regular data according to SendDataConfig settings"

# Row 28
$ws.Range("C28").Value = "SCC function is enabled"
$ws.Range("D28").Value = "Device error notification"
$ws.Range("E28").Value = "1. Generate device error (e.g., paper jam)`n2. Monitor Send Device Error communication"
$ws.Range("F28").Value = "MFP should notify server of device error with correct error code"

# Row 29
$ws.Range("C29").Value = "SCC function is enabled, error previously reported"
$ws.Range("D29").Value = "Error removal notification"
$ws.Range("E29").Value = "1. Generate device error`n2. Clear the error`n3. Monitor Send Device Error communication"
$ws.Range("F29").Value = "MFP should notify server of error removal with `"-`" prefix (e.g., -D102)"

# Row 30
$ws.Range("B30").Value = "TC025"
$ws.Range("C30").Value = "SCC function is enabled"
$ws.Range("D30").Value = "Panel message display"
$ws.Range("E30").Value = "1. Trigger SCC communication`n2. Observe panel display"
$ws.Range("F30").Value = "MFP should display appropriate message with status code during SCC processing"
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = ""

# Row 31
$ws.Range("B31").Value = "TC026"
$ws.Range("C31").Value = "SCC function is enabled, MFP in Super Sleep"
$ws.Range("D31").Value = "Super Sleep recovery for scheduled communication"
$ws.Range("E31").Value = "1. Put MFP in Super Sleep mode`n2. Wait for scheduled communication time"
$ws.Range("F31").Value = "MFP should wake from Super Sleep, perform communication, then return to Super Sleep"
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = ""

# Row 32
$ws.Range("B32").Value = "TC027"
$ws.Range("C32").Value = "SCC function is enabled"
$ws.Range("D32").Value = "Communication retry on error"
$ws.Range("E32").Value = "1. Disconnect network during communication`n2. Monitor retry behavior"
$ws.Range("F32").Value = "MFP should retry after 60 seconds, then enter retry mode if unsuccessful"
$ws.Range("G32").Value = ""
$ws.Range("H32").Value = ""

# Row 33
$ws.Range("B33").Value = "TC028"
$ws.Range("C33").Value = "SCC function is enabled, in retry mode"
$ws.Range("D33").Value = "Daily retry in retry mode"
$ws.Range("E33").Value = "1. Put MFP in retry mode`n2. Monitor communication attempts"
$ws.Range("F33").Value = "MFP should attempt communication once every 24 hours while in retry mode"
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = ""

# Row 34
$ws.Range("B34").Value = "TC029"
$ws.Range("C34").Value = "SCC function is enabled, custom URL configured"
$ws.Range("D34").Value = "First registration URL setting"
$ws.Range("E34").Value = "1. Set custom URL in 08-3827`n2. Enable SCC function`n3. Monitor registration process"
$ws.Range("F34").Value = "MFP should connect to specified custom URL instead of default server"
$ws.Range("G34").Value = ""
$ws.Range("H34").Value = ""

# Row 35
$ws.Range("B35").Value = "TC030"
$ws.Range("C35").Value = "SCC function is enabled, URL forward disabled"
$ws.Range("D35").Value = "URL forward setting"
$ws.Range("E35").Value = "1. Set 08-3828 to 0 (Disable)`n2. Configure server to return redirect URL`n3. Monitor registration process"
$ws.Range("F35").Value = "MFP should not follow redirect URL from server"
$ws.Range("G35").Value = ""
$ws.Range("H35").Value = ""

# Row 36
$ws.Range("B36").Value = "TC031"
$ws.Range("C36").Value = "SCC function is enabled, persistent policy configured"
$ws.Range("D36").Value = "Persistent policy check"
$ws.Range("E36").Value = "1. Configure persistent policy from server`n2. Change settings locally`n3. Wait for persistent check interval"
$ws.Range("F36").Value = "MFP should revert settings to match policy values"
$ws.Range("G36").Value = ""
$ws.Range("H36").Value = ""

# Row 37
$ws.Range("B37").Value = "TC032"
$ws.Range("C37").Value = "SCC function is enabled, time-based settings configured"
$ws.Range("D37").Value = "Time-based device state values"
$ws.Range("E37").Value = "1. Configure time-based settings from server`n2. Monitor settings at different times"
$ws.Range("F37").Value = "MFP should change settings according to configured time periods"
$ws.Range("G37").Value = ""
$ws.Range("H37").Value = ""

# Row 38
$ws.Range("B38").Value = "TC033"
$ws.Range("C38").Value = "SCC function is enabled"
$ws.Range("D38").Value = "IP redirect functionality"
$ws.Range("E38").Value = "1. Enable SCC on new MFP`n2. Monitor GetRedirectURL request`n3. Verify registration to appropriate regional server"
$ws.Range("F38").Value = "MFP should request redirect URL and connect to appropriate regional server"
$ws.Range("G38").Value = ""
$ws.Range("H38").Value = ""

# Row 39
$ws.Range("B39").Value = "TC034"
$ws.Range("C39").Value = "SCC function is enabled, communication in progress"
$ws.Range("D39").Value = "Power save prevention"
$ws.Range("E39").Value = "1. Trigger SCC communication`n2. Attempt to enter power save mode"
$ws.Range("F39").Value = "MFP should not enter Power Save, Sleep or Super Sleep during communication cycle"
$ws.Range("G39").Value = ""
$ws.Range("H39").Value = ""

# Row 40
$ws.Range("B40").Value = "TC035"
$ws.Range("C40").Value = "SCC function is enabled, communication in progress"
$ws.Range("D40").Value = "Service UI exclusion"
$ws.Range("E40").Value = "1. Start SCC communication`n2. Attempt to access service UI"
$ws.Range("F40").Value = "Service UI should be unavailable during SCC communication"
$ws.Range("G40").Value = ""
$ws.Range("H40").Value = ""

# Row 41
$ws.Range("B41").Value = "TC036"
$ws.Range("C41").Value = "SCC function is enabled, RDMS (eBR2) enabled"
$ws.Range("D41").Value = "Compatibility with RDMS"
$ws.Range("E41").Value = "1. Enable SCC function`n2. Enable RDMS (eBR2) function`n3. Test RDMS functionality"
$ws.Range("F41").Value = "RDMS should not work when SCC function is enabled (for versions before L6.02/L6.03)"
$ws.Range("G41").Value = ""
$ws.Range("H41").Value = ""

# Row 42
$ws.Range("B42").Value = "TC037"
$ws.Range("C42").Value = "SCC function is enabled, clone data instruction received"
$ws.Range("D42").Value = "Clone data application"
$ws.Range("E42").Value = "1. Configure server to send clone data instruction`n2. Trigger regular communication`n3. Monitor clone data application"
$ws.Range("F42").Value = "MFP should apply clone data regardless of SendDataConfig parameters"
$ws.Range("G42").Value = ""
$ws.Range("H42").Value = ""

# Row 43
$ws.Range("B43").Value = "TC038"
$ws.Range("C43").Value = "SCC function is enabled, invalid certificate"
$ws.Range("D43").Value = "HTTPS certificate validation"
$ws.Range("E43").Value = "1. Configure custom URL with invalid certificate`n2. Trigger SCC communication"
$ws.Range("F43").Value = "Communication should fail due to certificate validation"
$ws.Range("G43").Value = ""
$ws.Range("H43").Value = ""

# Row 44
$ws.Range("B44").Value = "TC039"
$ws.Range("C44").Value = "SCC function is enabled, HDD full"
$ws.Range("D44").Value = "HDD full handling"
$ws.Range("E44").Value = "1. Fill HDD to capacity`n2. Trigger package download"
$ws.Range("F44").Value = "MFP should delete downloaded data and exit communication cycle"
$ws.Range("G44").Value = ""
$ws.Range("H44").Value = ""

# Row 45
$ws.Range("B45").Value = "TC040"
$ws.Range("C45").Value = "SCC function is enabled, hash mismatch"
$ws.Range("D45").Value = "Hash validation"
$ws.Range("E45").Value = "1. Configure server to send package with incorrect hash`n2. Trigger download`n3. Monitor hash validation"
$ws.Range("F45").Value = "MFP should delete downloaded data and exit communication cycle"
$ws.Range("G45").Value = ""
$ws.Range("H45").Value = ""
